$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 450
$ws.Range("I2").Value = 433.66666
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 433.66666
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -320.66666
$ws.Range("N2").Value = -725

$ws.Range("H43").Value = 13075
$ws.Range("J43").Value = 683.1667
$ws.Range("L43").Value = 683.1667
$ws.Range("N43").Value = -821.1667

$ws.Range("H55").Value = 280
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 325
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 325
$ws.Range("M55").Value = -36
$ws.Range("N55").Value = -753

$ws.Range("H129").Value = 1142.1719
$ws.Range("I129").Value = 448.5
$ws.Range("J129").Value = 1164.5483
$ws.Range("K129").Value = 1345.5
$ws.Range("L129").Value = 3493.6449
$ws.Range("M129").Value = 3654.5
$ws.Range("N129").Value = -13493.6449

$ws.Range("H137").Value = 2523.6287
$ws.Range("I137").Value = 1557.1765
$ws.Range("J137").Value = 3436.389
$ws.Range("K137").Value = 4671.529500000001
$ws.Range("L137").Value = 10309.167
$ws.Range("M137").Value = -2121.529500000001
$ws.Range("N137").Value = -15409.167

$ws.Range("H138").Value = 3601.5796
$ws.Range("I138").Value = 2006.1154
$ws.Range("J138").Value = 4270.645
$ws.Range("K138").Value = 6018.3462
$ws.Range("L138").Value = 12811.935
$ws.Range("M138").Value = -878.3462
$ws.Range("N138").Value = -23091.935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2279.6667
$ws.Range("I88").Value = 1979.25
$ws.Range("K88").Value = 1979.25
$ws.Range("M88").Value = -1573.25

$ws.Range("H91").Value = 2279.6667
$ws.Range("I91").Value = 1979.25
$ws.Range("K91").Value = 1979.25
$ws.Range("M91").Value = -575.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70291.53
$ws.Range("I86").Value = 3586.875
$ws.Range("J86").Value = 146525.42
$ws.Range("K86").Value = 3586.875
$ws.Range("L86").Value = 146525.42
$ws.Range("M86").Value = -2463.875
$ws.Range("N86").Value = -148771.42

$ws.Range("H89").Value = 70291.53
$ws.Range("I89").Value = 3586.875
$ws.Range("J89").Value = 146525.42
$ws.Range("K89").Value = 17934.375
$ws.Range("L89").Value = 732627.1000000001
$ws.Range("M89").Value = -12318.375
$ws.Range("N89").Value = -743859.1000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 44758.75
$ws.Range("I62").Value = 53090.5
$ws.Range("J62").Value = 3100
$ws.Range("K62").Value = 53090.5
$ws.Range("L62").Value = 3100
$ws.Range("M62").Value = -52466.5
$ws.Range("N62").Value = -4348

$ws.Range("H65").Value = 44758.75
$ws.Range("I65").Value = 53090.5
$ws.Range("J65").Value = 3100
$ws.Range("K65").Value = 265452.5
$ws.Range("L65").Value = 15500
$ws.Range("M65").Value = -262332.5
$ws.Range("N65").Value = -21740

$ws.Range("H132").Value = 713866.5600000001
$ws.Range("I132").Value = 1230312.5
$ws.Range("J132").Value = 3753.5
$ws.Range("K132").Value = 3690937.5
$ws.Range("L132").Value = 11260.5
$ws.Range("M132").Value = -3688407.5
$ws.Range("N132").Value = -16320.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3815.6
$ws.Range("I116").Value = 359.33334
$ws.Range("K116").Value = 1078.00002
$ws.Range("M116").Value = 2363.99998

$ws.Range("H119").Value = 6499.8335
$ws.Range("J119").Value = 20000
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676

$ws.Range("H125").Value = 3523.75
$ws.Range("I125").Value = 300
$ws.Range("J125").Value = 3984.2856
$ws.Range("K125").Value = 900
$ws.Range("L125").Value = 11952.8568
$ws.Range("M125").Value = 4020
$ws.Range("N125").Value = -21792.8568

$ws.Range("H133").Value = 4812.6
$ws.Range("I133").Value = 1858.2
$ws.Range("J133").Value = 6289.8
$ws.Range("K133").Value = 5574.6
$ws.Range("L133").Value = 18869.4
$ws.Range("M133").Value = -514.6000000000004
$ws.Range("N133").Value = -28989.4

$ws.Range("H134").Value = 3378.5173
$ws.Range("I134").Value = 2339.8333
$ws.Range("J134").Value = 5078.1816
$ws.Range("K134").Value = 7019.499899999999
$ws.Range("L134").Value = 15234.5448
$ws.Range("M134").Value = -1949.499899999999
$ws.Range("N134").Value = -25374.5448

$ws.Range("H137").Value = 30307654
$ws.Range("I137").Value = 4447.75
$ws.Range("J137").Value = 47623772
$ws.Range("K137").Value = 13343.25
$ws.Range("L137").Value = 142871316
$ws.Range("M137").Value = -8243.25
$ws.Range("N137").Value = -142881516

$ws.Range("H139").Value = 2389.4
$ws.Range("I139").Value = 2236.75
$ws.Range("K139").Value = 6710.25
$ws.Range("M139").Value = -1570.25

$ws.Range("H140").Value = 1519
$ws.Range("I140").Value = 771.1739
$ws.Range("J140").Value = 3239
$ws.Range("K140").Value = 2313.5217
$ws.Range("L140").Value = 9717
$ws.Range("M140").Value = 2866.4783
$ws.Range("N140").Value = -20077

$ws.Range("H141").Value = 8757.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 18862.727
$ws.Range("J42").Value = 18862.727
$ws.Range("L42").Value = 18862.727
$ws.Range("N42").Value = -19832.727

$ws.Range("H70").Value = 309060.4
$ws.Range("I70").Value = 422587.47
$ws.Range("J70").Value = 6321.5557
$ws.Range("K70").Value = 422587.47
$ws.Range("L70").Value = 6321.5557
$ws.Range("M70").Value = -422317.47
$ws.Range("N70").Value = -6861.5557

$ws.Range("H73").Value = 309060.4
$ws.Range("I73").Value = 422587.47
$ws.Range("J73").Value = 6321.5557
$ws.Range("K73").Value = 422587.47
$ws.Range("L73").Value = 6321.5557
$ws.Range("M73").Value = -421651.47
$ws.Range("N73").Value = -8193.555700000001

$ws.Range("H102").Value = 3563.652
$ws.Range("I102").Value = 3545.9048
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 3545.9048
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -1923.9048
$ws.Range("N102").Value = -6994

$ws.Range("H115").Value = 18862.727
$ws.Range("J115").Value = 18862.727
$ws.Range("L115").Value = 18862.727
$ws.Range("N115").Value = -21212.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2706.75
$ws.Range("I7").Value = 2600
$ws.Range("J7").Value = 2742.3333
$ws.Range("K7").Value = 2600
$ws.Range("L7").Value = 2742.3333
$ws.Range("M7").Value = -2488
$ws.Range("N7").Value = -2966.3333

$ws.Range("H42").Value = 34248.75
$ws.Range("I42").Value = 38500
$ws.Range("J42").Value = 29997.5
$ws.Range("K42").Value = 38500
$ws.Range("L42").Value = 29997.5
$ws.Range("M42").Value = -37937
$ws.Range("N42").Value = -31123.5

$ws.Range("H43").Value = 45000
$ws.Range("J43").Value = 40000
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40386

$ws.Range("H49").Value = 34248.75
$ws.Range("I49").Value = 38500
$ws.Range("J49").Value = 29997.5
$ws.Range("K49").Value = 38500
$ws.Range("L49").Value = 29997.5
$ws.Range("M49").Value = -38353
$ws.Range("N49").Value = -30291.5

$ws.Range("H115").Value = 96100.664
$ws.Range("J115").Value = 96100.664
$ws.Range("L115").Value = 96100.664
$ws.Range("N115").Value = -98450.664

$ws.Range("H126").Value = 2706.75
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 2742.3333
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 8226.999899999999
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -13166.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 104168616
$ws.Range("I122").Value = 125001870
$ws.Range("J122").Value = 2327.5
$ws.Range("K122").Value = 375005610
$ws.Range("L122").Value = 6982.5
$ws.Range("M122").Value = -375003160
$ws.Range("N122").Value = -11882.5
